$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 12 with trade data (mirrors the row 11 formatting, esp. the
# date-formatted column G)
$ws.Range("A12").Value = 10279.24
$ws.Range("B12").Value = 10189.57
$ws.Range("C12").Value = 105.78
$ws.Range("D12").Value = 106.71
$ws.Range("E12").Value = $false
$ws.Range("F12").Value = 0.88

# Copy formatting from G11 (date number format) before setting the new value
$ws.Range("G11").Copy()
$ws.Range("G12").PasteSpecial(-4122)
$ws.Range("G12").Value = 42620.766238425924

$ws.Range("H12").Value = $true
